$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (date) and C (id) hold text-like values ("2020-02-27", "0216")
# that Excel's normal type-inference would coerce into a date serial / a
# number. Force text interpretation via NumberFormat "@" before assigning,
# then restore the default "Normal" style so the new rows end up with no
# explicit style attribute, matching the rest of the sheet.
$ws.Range("B69:C70").NumberFormat = "@"

$ws.Range("A69").Value = 1582761600
$ws.Range("B69").Value = "2020-02-27"
$ws.Range("C69").Value = "0216"
$ws.Range("D69").Value = "SPRING"
$ws.Range("E69").Value = 0.23
$ws.Range("F69").Value = 0.23
$ws.Range("G69").Value = 0.225
$ws.Range("H69").Value = 0.23
$ws.Range("I69").Value = 209600

$ws.Range("A70").Value = 1582848000
$ws.Range("B70").Value = "2020-02-28"
$ws.Range("C70").Value = "0216"
$ws.Range("D70").Value = "SPRING"
$ws.Range("E70").Value = 0.225
$ws.Range("F70").Value = 0.235
$ws.Range("G70").Value = 0.22
$ws.Range("H70").Value = 0.225
$ws.Range("I70").Value = 1385800

$ws.Range("B69:C70").Style = "Normal"
